$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the dataset. It is inserted
# as a new row 603, pushing all the former rows 603-672 down by one
# (to 604-673), exactly like Excel's normal "insert row" behaviour which
# also grows the sheet's used range from A1:R672 to A1:R673.
$ws.Rows("603:603").Insert()

# Populate the newly inserted row 603 with the new record's data.
$ws.Range("A603").Value2 = 5
$ws.Range("B603").Value2 = 'Macroferia Regional de Talca'
$ws.Range("C603").Value2 = 'Maule'
$ws.Range("D603").Value2 = 45212
$ws.Range("E603").Value2 = 7
$ws.Range("F603").Value2 = 100114014
$ws.Range("G603").Value2 = 'Betarraga'
$ws.Range("H603").Value2 = 'Sin especificar'
$ws.Range("I603").Value2 = 'Primera'
$ws.Range("J603").Value2 = 4000
$ws.Range("K603").Value2 = 500
$ws.Range("L603").Value2 = 500
$ws.Range("M603").Value2 = 500
$ws.Range("N603").Value2 = '$/paquete 5 unidades'
$ws.Range("O603").Value2 = 'Región del Maule'
$ws.Range("P603").Value2 = 100
$ws.Range("Q603").Value2 = 5
$ws.Range("R603").Value2 = 'Hortaliza'
